$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 4-9: column A becomes "Course N" (N = row-3),
# column B alternates VO (odd course numbers) / VU (even course numbers).
$ws.Range("B4").Value = "VO"
$ws.Range("B5").Value = "VU"
$ws.Range("B6").Value = "VO"
$ws.Range("B8").Value = "VO"
$ws.Range("A9").Value = "Course 6"
$ws.Range("B9").Value = "VU"

# Append new rows 10-13, copying the number format of column H from row 9
# so the new date cells keep the same built-in short-date style.
$rows = @(
    @{ Row = 10; Course = "Course 7";  B = "VO"; D = 7;  H = 43471 },
    @{ Row = 11; Course = "Course 8";  B = "VU"; D = 8;  H = 43471 },
    @{ Row = 12; Course = "Course 9";  B = "VO"; D = 9;  H = 43471 },
    @{ Row = 13; Course = "Course 10"; B = "VU"; D = 10; H = 43471 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Course
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = "A"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = "sehr gut"
    $ws.Range("F$row").Value = 1
    $ws.Range("G$row").Value = 1

    $ws.Range("H9").Copy()
    $ws.Range("H$row").PasteSpecial(-4122)
    $ws.Range("H$row").Value = $r.H

    $ws.Range("I$row").Value = "Person"
}

$ws.Range("B13").Select()
